$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = "Objects, Variables & Operators -- HP"
$ws.Range("E4").Value = "Control Structures: Conditional Statements -- HP"
$ws.Range("E5").Value = "Control Structures: Loops -- HP"
$ws.Range("E6").Value = "Consolidation 1 - Turtlesim? Text-based adventure game? + debugging -- HP"
$ws.Range("E7").Value = "Data Structures (Lists) -- MH"
$ws.Range("E9").Value = "Functions 1 (define, returning) -- MH"
$ws.Range("E10").Value = "Reading and plotting data, Matplotlib - MH"
$ws.Range("E11").Value = "Consolidation 2 - descriptive statistics of a data set (w/ independent use of a python module for statistical analysis) - MH"
$ws.Range("E12").Value = "Functions 2 (scope) -- MH"
$ws.Range("E13").Value = "Coursework support session  - HP + MH"
$ws.Range("E14").Value = "Consolidation 3 - Refactoring and extending code -> Testing and ChatGPT - HP"

$ws.Range("E15").Select()
